$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = -10.85
$ws.Range("C18").Value = -11.62
$ws.Range("C20").Value = -12.354
$ws.Range("C27").Value = -12.951
$ws.Range("C69").Value = -10.537
$ws.Range("C76").Value = -13.032
$ws.Range("C82").Value = -11.885
